# Sprint 4.xlsx update — team re-assignment on the "Shopenzer Testcases" sheet.
# The team ID and the four team-member names (Executed By column) were
# refreshed, and the workbook was left with the testcases sheet active /
# selected at N9 (the last "Executed By" cell), matching how the sheet was
# left after the edit.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Shopenzer Testcases")
$ws2 = $wb.Worksheets.Item("Testscearnios")

# --- Update the Team ID (row 2, column F on the summary block) ---
$ws1.Range("F2").Value = "PNT2022TMID53380"

# --- Update the "Executed By" names for the four sprint test cases ---
$ws1.Range("N6").Value = "Ritunjay M"
$ws1.Range("N7").Value = "Praveen Raagul R"
$ws1.Range("N8").Value = "Pradeep V"
$ws1.Range("N9").Value = "Munish Kumar S"

# --- Leave the workbook with "Shopenzer Testcases" active and N9 selected ---
$ws1.Activate()
$ws1.Range("N9").Select()

# Zoom tweak captured alongside the save (matches the refreshed view state).
$win = $excel.ActiveWindow
$win.Zoom = 57
